$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 77

function Set-TextCell($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

Set-TextCell $ws.Cells.Item($row, 1) "2023-06-25"
Set-TextCell $ws.Cells.Item($row, 2) "22:53:45"
Set-TextCell $ws.Cells.Item($row, 3) "Sunday"
Set-TextCell $ws.Cells.Item($row, 4) "26"

$ws.Cells.Item($row, 5).Value = 122704
$ws.Cells.Item($row, 6).Value = 134220
$ws.Cells.Item($row, 7).Value = 163510
$ws.Cells.Item($row, 8).Value = 133592
$ws.Cells.Item($row, 9).Value = 177287
$ws.Cells.Item($row, 10).Value = 116615
$ws.Cells.Item($row, 11).Value = 203560
$ws.Cells.Item($row, 12).Value = 226246
$ws.Cells.Item($row, 13).Value = 175924
$ws.Cells.Item($row, 14).Value = 104241
$ws.Cells.Item($row, 15).Value = 39650
$ws.Cells.Item($row, 16).Value = 33823
$ws.Cells.Item($row, 17).Value = 52124
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 36322
$ws.Cells.Item($row, 20).Value = -1
